# Update Active_Outages.xlsx - 6/19/2025, 10:14:45 AM
#
# This refreshes the "Elapsed Duration(Hrs)" counters on several open
# outage tickets, drops the now-resolved spacer/blank row from the R1
# ("R1") sheet, and appends the ticket that moved from R1 into R2
# ("R2") as its new last row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Refresh "Elapsed Duration(Hrs)" values (column G) across sheets
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3947:28:51"
$ws.Range("G3").Value = "87:01:29"
$ws.Range("G4").Value = "110:01:29"

$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12128:52:32"
$ws.Range("G3").Value = "3258:36:01"
$ws.Range("G4").Value = "496:47:35"

$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2974:42:21"
$ws.Range("G3").Value = "201:54:36"
$ws.Range("G4").Value = "90:07:01"
$ws.Range("G5").Value = "87:44:34"

$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "448:41:20"

$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "89:13:38"

# ---------------------------------------------------------------
# 2) R1 sheet ("R1"): remove the trailing blank row 6
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Rows("6").Delete()

# ---------------------------------------------------------------
# 3) R2 sheet ("R2"): append the new outage row (row 6)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("B6").Value = "R4"
$ws2.Range("D6").Value = "JED0125"
$ws2.Range("I6").Value = "Generator-SG"
$ws2.Range("J6").Value = "Good+In progress"
$ws2.Range("L6").Value = "Latis"
